$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.312.60'
$ws.Range("E2").Value = '  +1.83%  '

$ws.Range("D3").Value = '3.596.14'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '656.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '

$ws.Range("E7").Value = '  +12.79%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.414'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.06'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = '3.592.44'
$ws.Range("E11").Value = '  -0.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.62%  '

$ws.Range("E13").Value = '  +1.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '

$ws.Range("D15").Value = '4.264.95'
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("D16").Value = '97.090.87'
$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("E17").Value = '  +1.99%  '

$ws.Range("D18").Value = '3.587.26'
$ws.Range("E18").Value = '  -0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.529'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '511.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.38%  '

$ws.Range("E24").Value = '  -3.30%  '

$ws.Range("E25").Value = '  +2.35%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '97.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.31%  '

$ws.Range("D29").Value = '3.789.01'
$ws.Range("E29").Value = '  -0.83%  '

$ws.Range("E30").Value = '  -1.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.151'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.75%  '

$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("E34").Value = '  +4.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.65'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.82%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '624.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.83%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.86%  '

$ws.Range("E39").Value = '  +2.19%  '

$ws.Range("E40").Value = '  +11.21%  '

$ws.Range("E41").Value = '  +1.44%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.916'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("E44").Value = '  +6.11%  '

$ws.Range("E45").Value = '  +5.12%  '

$ws.Range("E46").Value = '  +3.15%  '

$ws.Range("E47").Value = '  +1.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.06%  '

$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.96%  '
